$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Text-run normalization: re-run a self Find/Replace across the exact
#    span of adjacent, identically-formatted runs so Word recombines them
#    into a single run (no visible text change, only run-splitting changes).
# ---------------------------------------------------------------------------
function Merge-Runs([string]$text) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $rng.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, $text, 2) | Out-Null
}

Merge-Runs "              To write a C program to simulate a Deterministic Finite Automata (DFA) for the given language representing strings that start with a and end with "
Merge-Runs "   1. Draw a DFA for the given language and construct the transition table."
Merge-Runs "   2. Store the transition table in a two-dimensional array. "
Merge-Runs "   3. Initialize "
Merge-Runs "   4. Get the input string from the user."
Merge-Runs "   5. Find the length of the input string. "
Merge-Runs "   6. Read the input string character by character. "
Merge-Runs "   7. Repeat step 8 for every character"
Merge-Runs "   8. Refer the transition table for the entry corresponding to the present state and the current input symbol and update the next state."
Merge-Runs "   9. When we reach the end of the input, if the final state is reached, the input is accepted. "
Merge-Runs " #include<stdio.h>"

# ---------------------------------------------------------------------------
# 2. Picture 1 (Output screenshot #1): crop top + shrink height, drop border
# ---------------------------------------------------------------------------
$pic1 = $d.InlineShapes.Item(1)
$pic1.PictureFormat.CropTop = 25.25112
$pic1.LockAspectRatio = 0
$pic1.Width = 409.7999212598425
$pic1.Height = 103.2
$pic1.LockAspectRatio = 1
$pic1.Line.Visible = $false

# ---------------------------------------------------------------------------
# 3. Picture 2 (Output screenshot #2): crop top + shrink height, drop border
# ---------------------------------------------------------------------------
$pic2 = $d.InlineShapes.Item(2)
$pic2.PictureFormat.CropTop = 27.528495
$pic2.LockAspectRatio = 0
$pic2.Width = 451.3
$pic2.Height = 111.4
$pic2.LockAspectRatio = 1
$pic2.Line.Visible = $false
